$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsEoC   = $wb.Worksheets.Item("EoCEDwEC")

# --- EoCEDwEC sheet: add the four new fuel rows (string table order matters) ---
$wsEoC.Range("A8").Value = "kerosene"
$wsEoC.Range("A9").Value = "heavy or residual fuel oil"
$wsEoC.Range("A10").Value = "LPG propane or butane"
$wsEoC.Range("A11").Value = "hydrogen"

# --- About sheet: add explanatory notes about new fuel types (rows 28-31) ---
$wsAbout.Range("A28").Value = "We assume kerosene and fuel oil-burning equipment is similar to diesel-burning equipment."
$wsAbout.Range("A29").Value = "We assume LPG/propane/butane-burning equipment is similar to natural gas-burning equipment."
$wsAbout.Range("A30").Value = "We assume hydrogen-using equipment is similar to electricity-using equipment (as it may"
$wsAbout.Range("A31").Value = "contain fuel cells that produce electricity from hydrogen)."

# --- EoCEDwEC sheet: rename header ---
$wsEoC.Range("A1").Value = "Elasticity by Fuel (dimensionless)"
$wsEoC.Range("A1").Font.Bold = $true
$wsEoC.Range("A1").WrapText = $true
$wsEoC.Range("A1").EntireRow.RowHeight = 30

# Match the orange fill used by the "coal" row (row 3), which is also a
# fuel whose elasticity is derived from another fuel's values.
$wsEoC.Range("B8:D11").Interior.Color = $wsEoC.Range("B3").Interior.Color

# kerosene & heavy/residual fuel oil -> same as petroleum diesel (row 5)
$wsEoC.Range("B8").Formula = "=B5"
$wsEoC.Range("C8").Formula = "=C5"
$wsEoC.Range("D8").Formula = "=D5"

$wsEoC.Range("B9").Formula = "=B5"
$wsEoC.Range("C9").Formula = "=C5"
$wsEoC.Range("D9").Formula = "=D5"

# LPG/propane/butane -> same as natural gas (row 4)
$wsEoC.Range("B10").Formula = "=B4"
$wsEoC.Range("C10").Formula = "=C4"
$wsEoC.Range("D10").Formula = "=D4"

# hydrogen -> same as electricity (row 2)
$wsEoC.Range("B11").Formula = "=B2"
$wsEoC.Range("C11").Formula = "=C2"
$wsEoC.Range("D11").Formula = "=D2"

$wb.Save()
